$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '65.396.10'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.25%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.388.22'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.71%  '

# Row 4
$ws.Range('E4').Value = '  +0.14%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '558.63'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.47%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '175.46'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.60%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.631'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.01%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.376.29'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.50%  '

# Row 9
$ws.Range('E9').Value = '  +0.21%  '

# Row 10
$ws.Range('E10').Value = '  +4.53%  '

# Row 11
$ws.Range('E11').Value = '  +0.75%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '53.30'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -4.38%  '

# Row 13
$ws.Range('E13').Value = '  +0.16%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.20'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.74%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.929.48'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.08%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '18.29'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.24%  '

# Row 17
$ws.Range('E17').Value = '  +1.34%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.378.60'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.62%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '65.235.38'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.23%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.80'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.41%  '

# Row 21
$ws.Range('E21').Value = '  +0.73%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '488.44'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +5.52%  '

# Row 23
$ws.Range('E23').Value = '  -0.23%  '

# Row 24
$ws.Range('E24').Value = '  +0.29%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '14.23'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +4.78%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '87.87'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.85%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.90'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.82%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.69'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.83%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.70'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.46%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '31.27'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.41%  '

# Row 31
$ws.Range('E31').Value = '  -2.53%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '11.46'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.48%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '62.72'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.73%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '574.22'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.35%  '

# Row 35
$ws.Range('E35').Value = '  -0.42%  '

# Row 36
$ws.Range('E36').Value = '  -0.08%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.60'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.64%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.140'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.08%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '35.80'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.74%  '

# Row 40
$ws.Range('E40').Value = '  -0.04%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0₃0737'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.60%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.124.57'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.08%  '

# Row 43
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0417'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.88%  '

# Row 44
$ws.Range('B44').Value = 'ThetaToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.78'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.77%  '

# Row 45
$ws.Range('E45').Value = '  +1.04%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.17'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.76%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.43'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.97%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.999'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.20%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '140.02'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.47%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.56'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.23%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.41'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.12%  '
